# Mise à jour du rapport de projet :
# - restaurer la fenêtre du classeur (elle n'est plus réduite/minimisée)
# - rafraîchir la vue de la feuille PlanningProjet (zoom, défilement, sélection)
# - avancer le "Semaine_Affichage" d'une semaine (C4 : -2 -> -3)
# - renseigner la date de fin de la tâche "Création du document rapport de projet" (F25)
# NB : C3 contient =TODAY(), sa valeur se recalcule automatiquement.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PlanningProjet")

# Restaurer la fenêtre (elle était minimisée)
$win = $excel.ActiveWindow
$win.WindowState = -4143   # xlNormal

# Nouveau zoom de la feuille (70% -> 85%)
$win.Zoom = 85

# Décaler l'affichage de la "semaine" de -2 à -3
$ws.Range("C4").Value = -3

# Renseigner la date de fin manquante pour la ligne 25
$ws.Range("F25").Value = 44705

# Repositionner la sélection / le défilement sur la feuille
$ws.Range("A16").Select()
$win.ScrollRow = 16
$ws.Range("F25").Select()
